# Update the "STEM" label in the research hours table to "S-STEM"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("A2").Value = "S-STEM"

# Move the active selection to F4, matching the saved cursor position
$ws.Range("F4").Select()
